$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").ClearContents()

$ws.Range("B4").Value = "Philipot"
$ws.Range("D4").Value = "Philipot"
$ws.Range("E4").Value = "Montebourg"
$ws.Range("F4").Value = "Montebourg"
$ws.Range("G4").Value = "Montebourg"
$ws.Range("H4").Value = "Philipot"
$ws.Range("I4").Value = "Montebourg"
$ws.Range("J4").Value = "Philipot"
$ws.Range("K4").Value = "Montebourg"
$ws.Range("L4").Value = "Montebourg"
$ws.Range("M4").Value = "Montebourg"
$ws.Range("N4").Value = "Philipot"
$ws.Range("O4").Value = "Montebourg"
$ws.Range("P4").Value = "Philipot"
$ws.Range("Q4").Value = "Montebourg"
$ws.Range("R4").Value = "Montebourg"
$ws.Range("S4").Value = "Montebourg"
$ws.Range("T4").Value = "Macron"
$ws.Range("U4").Value = "Montebourg"
$ws.Range("V4").Value = "Philipot"
$ws.Range("W4").Value = "Montebourg"
$ws.Range("X4").Value = "Montebourg"
$ws.Range("Y4").Value = "Montebourg"
$ws.Range("Z4").Value = "Macron"
$ws.Range("AB4").Value = "Philipot"
$ws.Range("AF4").Value = "Macron"
$ws.Range("AH4").Value = "Philipot"
$ws.Range("AN4").Value = "Philipot"
$ws.Range("B5").Value = "Philipot"
$ws.Range("C5").Value = "Philipot"
$ws.Range("D5").Value = "Philipot"
$ws.Range("E5").Value = "Montebourg"
$ws.Range("F5").Value = "Montebourg"
$ws.Range("G5").Value = "Montebourg"
$ws.Range("H5").Value = "Philipot"
$ws.Range("I5").Value = "Montebourg"
$ws.Range("J5").Value = "Montebourg"
$ws.Range("K5").Value = "Montebourg"
$ws.Range("L5").Value = "Montebourg"
$ws.Range("M5").Value = "Montebourg"
$ws.Range("N5").Value = "Philipot"
$ws.Range("O5").Value = "Philipot"
$ws.Range("P5").Value = "Philipot"
$ws.Range("Q5").Value = "Philipot"
$ws.Range("R5").Value = "Philipot"
$ws.Range("S5").Value = "Philipot"
$ws.Range("T5").Value = "Macron"
$ws.Range("U5").Value = "Macron"
$ws.Range("V5").Value = "Macron"
$ws.Range("W5").Value = "Macron"
$ws.Range("X5").Value = "Macron"
$ws.Range("Y5").Value = "Macron"
$ws.Range("Z5").Value = "Macron"
$ws.Range("AA5").Value = "Macron"
$ws.Range("AB5").Value = "Macron"
$ws.Range("AC5").Value = "Macron"
$ws.Range("AD5").Value = "Macron"
$ws.Range("AE5").Value = "Macron"
$ws.Range("AF5").Value = "Macron"
$ws.Range("AG5").Value = "Macron"
$ws.Range("AH5").Value = "Macron"
$ws.Range("AI5").Value = "Macron"
$ws.Range("AJ5").Value = "Macron"
$ws.Range("AK5").Value = "Macron"
$ws.Range("B6").Value = "Philipot"
$ws.Range("C6").Value = "Philipot"
$ws.Range("D6").Value = "Philipot"
$ws.Range("E6").Value = "Montebourg"
$ws.Range("F6").Value = "Montebourg"
$ws.Range("G6").Value = "Montebourg"
$ws.Range("H6").Value = "Philipot"
$ws.Range("I6").Value = "Philipot"
$ws.Range("J6").Value = "Philipot"
$ws.Range("K6").Value = "Montebourg"
$ws.Range("L6").Value = "Montebourg"
$ws.Range("M6").Value = "Montebourg"
$ws.Range("N6").Value = "Philipot"
$ws.Range("O6").Value = "Philipot"
$ws.Range("P6").Value = "Philipot"
$ws.Range("Q6").Value = "Montebourg"
$ws.Range("R6").Value = "Montebourg"
$ws.Range("S6").Value = "Philipot"
$ws.Range("T6").Value = "Macron"
$ws.Range("U6").Value = "Montebourg"
$ws.Range("V6").Value = "Montebourg"
$ws.Range("W6").Value = "Montebourg"
$ws.Range("X6").Value = "Montebourg"
$ws.Range("Y6").Value = "Macron"
$ws.Range("Z6").Value = "Macron"
$ws.Range("AC6").Value = "Montebourg"
$ws.Range("AD6").Value = "Montebourg"
$ws.Range("AE6").Value = "Macron"
$ws.Range("AF6").Value = "Macron"
$ws.Range("AG6").Value = "Montebourg"
$ws.Range("AH6").Value = "Montebourg"
$ws.Range("AI6").Value = "Montebourg"
$ws.Range("AJ6").Value = "Montebourg"
$ws.Range("AK6").Value = "Macron"
$ws.Range("B7").Value = "Philipot"
$ws.Range("C7").Value = "Philipot"
$ws.Range("D7").Value = "Montebourg"
$ws.Range("E7").Value = "Montebourg"
$ws.Range("F7").Value = "Montebourg"
$ws.Range("G7").Value = "Montebourg"
$ws.Range("H7").Value = "Philipot"
$ws.Range("I7").Value = "Montebourg"
$ws.Range("J7").Value = "Montebourg"
$ws.Range("K7").Value = "Montebourg"
$ws.Range("L7").Value = "Montebourg"
$ws.Range("M7").Value = "Montebourg"
$ws.Range("N7").Value = "Montebourg"
$ws.Range("O7").Value = "Montebourg"
$ws.Range("P7").Value = "Montebourg"
$ws.Range("Q7").Value = "Montebourg"
$ws.Range("R7").Value = "Montebourg"
$ws.Range("S7").Value = "Montebourg"
$ws.Range("T7").Value = "Philipot"
$ws.Range("U7").Value = "Montebourg"
$ws.Range("V7").Value = "Montebourg"
$ws.Range("W7").Value = "Montebourg"
$ws.Range("X7").Value = "Montebourg"
$ws.Range("Y7").Value = "Montebourg"
$ws.Range("Z7").Value = "Montebourg"
$ws.Range("AF7").Value = "Hidalgo"
$ws.Range("AL7").Value = "Macron"
$ws.Range("AM7").Value = "Pécresse"
$ws.Range("AN7").Value = "Macron"
$ws.Range("B8").Value = "Philipot"
$ws.Range("C8").Value = "Philipot"
$ws.Range("D8").Value = "Philipot"
$ws.Range("E8").Value = "Montebourg"
$ws.Range("F8").Value = "Montebourg"
$ws.Range("G8").Value = "Montebourg"
$ws.Range("H8").Value = "Philipot"
$ws.Range("I8").Value = "Montebourg"
$ws.Range("J8").Value = "Montebourg"
$ws.Range("K8").Value = "Montebourg"
$ws.Range("L8").Value = "Montebourg"
$ws.Range("M8").Value = "Montebourg"
$ws.Range("N8").Value = "Montebourg"
$ws.Range("O8").Value = "Montebourg"
$ws.Range("P8").Value = "Montebourg"
$ws.Range("Q8").Value = "Montebourg"
$ws.Range("R8").Value = "Montebourg"
$ws.Range("S8").Value = "Montebourg"
$ws.Range("T8").Value = "Montebourg"
$ws.Range("U8").Value = "Montebourg"
$ws.Range("V8").Value = "Montebourg"
$ws.Range("W8").Value = "Montebourg"
$ws.Range("X8").Value = "Montebourg"
$ws.Range("Y8").Value = "Montebourg"
$ws.Range("Z8").Value = "Montebourg"
$ws.Range("AF8").Value = "Montebourg"
$ws.Range("AL8").Value = "Montebourg"
$ws.Range("B9").Value = "Philipot"
$ws.Range("C9").Value = "Philipot"
$ws.Range("D9").Value = "Philipot"
$ws.Range("E9").Value = "Montebourg"
$ws.Range("F9").Value = "Montebourg"
$ws.Range("G9").Value = "Montebourg"
$ws.Range("H9").Value = "Philipot"
$ws.Range("I9").Value = "Philipot"
$ws.Range("J9").Value = "Philipot"
$ws.Range("K9").Value = "Montebourg"
$ws.Range("L9").Value = "Montebourg"
$ws.Range("M9").Value = "Montebourg"
$ws.Range("N9").Value = "Philipot"
$ws.Range("O9").Value = "Philipot"
$ws.Range("P9").Value = "Philipot"
$ws.Range("Q9").Value = "Montebourg"
$ws.Range("R9").Value = "Montebourg"
$ws.Range("S9").Value = "Montebourg"
$ws.Range("T9").Value = "Macron"
$ws.Range("U9").Value = "Montebourg"
$ws.Range("V9").Value = "Montebourg"
$ws.Range("W9").Value = "Montebourg"
$ws.Range("X9").Value = "Montebourg"
$ws.Range("Y9").Value = "Montebourg"
$ws.Range("Z9").Value = "Macron"
$ws.Range("AF9").Value = "Macron"
$ws.Range("AK9").Value = "Macron"
$ws.Range("AM9").Value = "Macron"
$ws.Range("B10").Value = "Philipot"
$ws.Range("C10").Value = "Philipot"
$ws.Range("D10").Value = "Philipot"
$ws.Range("E10").Value = "Montebourg"
$ws.Range("F10").Value = "Montebourg"
$ws.Range("G10").Value = "Montebourg"
$ws.Range("H10").Value = "Philipot"
$ws.Range("I10").Value = "Philipot"
$ws.Range("J10").Value = "Philipot"
$ws.Range("K10").Value = "Montebourg"
$ws.Range("L10").Value = "Montebourg"
$ws.Range("M10").Value = "Montebourg"
$ws.Range("N10").Value = "Philipot"
$ws.Range("O10").Value = "Philipot"
$ws.Range("P10").Value = "Philipot"
$ws.Range("Q10").Value = "Montebourg"
$ws.Range("R10").Value = "Montebourg"
$ws.Range("S10").Value = "Montebourg"
$ws.Range("T10").Value = "Philipot"
$ws.Range("U10").Value = "Philipot"
$ws.Range("V10").Value = "Philipot"
$ws.Range("W10").Value = "Montebourg"
$ws.Range("X10").Value = "Montebourg"
$ws.Range("Y10").Value = "Montebourg"
$ws.Range("Z10").Value = "Philipot"
$ws.Range("AA10").Value = "Philipot"
$ws.Range("AB10").Value = "Philipot"
$ws.Range("AG10").Value = "Philipot"
$ws.Range("AH10").Value = "Philipot"
$ws.Range("AM10").Value = "Philipot"
$ws.Range("AN10").Value = "Philipot"
$ws.Range("B11").Value = "Philipot"
$ws.Range("C11").Value = "Montebourg"
$ws.Range("D11").Value = "Montebourg"
$ws.Range("E11").Value = "Montebourg"
$ws.Range("F11").Value = "Montebourg"
$ws.Range("G11").Value = "Montebourg"
$ws.Range("H11").Value = "Philipot"
$ws.Range("I11").Value = "Montebourg"
$ws.Range("J11").Value = "Montebourg"
$ws.Range("K11").Value = "Montebourg"
$ws.Range("L11").Value = "Montebourg"
$ws.Range("M11").Value = "Montebourg"
$ws.Range("N11").Value = "Philipot"
$ws.Range("O11").Value = "Montebourg"
$ws.Range("P11").Value = "Montebourg"
$ws.Range("Q11").Value = "Montebourg"
$ws.Range("R11").Value = "Montebourg"
$ws.Range("S11").Value = "Montebourg"
$ws.Range("T11").Value = "Macron"
$ws.Range("U11").Value = "Montebourg"
$ws.Range("V11").Value = "Montebourg"
$ws.Range("W11").Value = "Montebourg"
$ws.Range("X11").Value = "Montebourg"
$ws.Range("Y11").Value = "Montebourg"
$ws.Range("Z11").Value = "Macron"
$ws.Range("AF11").Value = "Macron"
$ws.Range("AL11").Value = "Macron"
$ws.Range("B12").Value = "Philipot"
$ws.Range("C12").Value = "Philipot"
$ws.Range("D12").Value = "Philipot"
$ws.Range("E12").Value = "Montebourg"
$ws.Range("F12").Value = "Montebourg"
$ws.Range("G12").Value = "Montebourg"
$ws.Range("H12").Value = "Philipot"
$ws.Range("I12").Value = "Montebourg"
$ws.Range("J12").Value = "Montebourg"
$ws.Range("K12").Value = "Montebourg"
$ws.Range("L12").Value = "Montebourg"
$ws.Range("M12").Value = "Montebourg"
$ws.Range("N12").Value = "Montebourg"
$ws.Range("O12").Value = "Montebourg"
$ws.Range("P12").Value = "Montebourg"
$ws.Range("Q12").Value = "Montebourg"
$ws.Range("R12").Value = "Montebourg"
$ws.Range("S12").Value = "Montebourg"
$ws.Range("T12").Value = "Montebourg"
$ws.Range("U12").Value = "Montebourg"
$ws.Range("V12").Value = "Montebourg"
$ws.Range("W12").Value = "Montebourg"
$ws.Range("X12").Value = "Montebourg"
$ws.Range("Y12").Value = "Montebourg"
$ws.Range("Z12").Value = "Montebourg"
$ws.Range("AF12").Value = "Montebourg"
$ws.Range("AL12").Value = "Montebourg"
